# New weekly price observation for "Haba" at Vega Monumental Concepción.
# A row is inserted above the current row 58, shifting the existing
# rows 58-61 down to 59-62 (and extending the used range to A1:R62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(58).Insert()

$ws.Range("A58").Value = 11
$ws.Range("B58").Value = "Vega Monumental Concepción"
$ws.Range("C58").Value = "Bíobío"
$ws.Range("D58").Value = 45147
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 100112026
$ws.Range("G58").Value = "Haba"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 50
$ws.Range("K58").Value = 14000
$ws.Range("L58").Value = 14000
$ws.Range("M58").Value = 14000
$ws.Range("N58").Value = '$/saco 25 kilos'
$ws.Range("O58").Value = "Región de Coquimbo"
$ws.Range("P58").Value = 560
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
